$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store plain text values (e.g. "295.41",
# "1.63%") rather than numbers, matching the workbook's original inline-string
# cell format. Temporarily force a text number format on the affected range so
# that assigning these numeric-looking / percent-looking strings does not get
# auto-converted by Excel into actual numeric/percentage values.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "295.41"
$ws.Range("E2").Value = "1.63%"
$ws.Range("D3").Value = "31.10"
$ws.Range("E3").Value = "0.66%"
$ws.Range("D4").Value = "4.933"
$ws.Range("E4").Value = "-0.30%"
$ws.Range("D5").Value = "0.07448"
$ws.Range("D6").Value = "2.244"
$ws.Range("E6").Value = "24.65%"
$ws.Range("D7").Value = "7.766"
$ws.Range("E7").Value = "1.35%"
$ws.Range("D8").Value = "3.744"
$ws.Range("E8").Value = "0.33%"
$ws.Range("D9").Value = "0.9134"
$ws.Range("E9").Value = "2.07%"
$ws.Range("D10").Value = "0.09058"
$ws.Range("E10").Value = "17.84%"
$ws.Range("D11").Value = "0.1715"
$ws.Range("E11").Value = "3.87%"
$ws.Range("D12").Value = "0.08283"
$ws.Range("E12").Value = "2.12%"
$ws.Range("D13").Value = "0.03124"
$ws.Range("E13").Value = "3.23%"
$ws.Range("E14").Value = "0.31%"
$ws.Range("E15").Value = "1.84%"
$ws.Range("D16").Value = "0.005698"
$ws.Range("E16").Value = "-1.22%"
$ws.Range("D17").Value = "3.499"
$ws.Range("E17").Value = "0.57%"
$ws.Range("D18").Value = "2.075"
$ws.Range("E18").Value = "-0.35%"
$ws.Range("D19").Value = "0.3328"
$ws.Range("E19").Value = "1.51%"
$ws.Range("D20").Value = "0.1297"
$ws.Range("E20").Value = "1.95%"
$ws.Range("D21").Value = "3.982"
$ws.Range("E21").Value = "-1.34%"
$ws.Range("E22").Value = "5.07%"
$ws.Range("D23").Value = "0.04545"
$ws.Range("E23").Value = "0.82%"
$ws.Range("E24").Value = "0.07%"
$ws.Range("D25").Value = "0.004616"
$ws.Range("E25").Value = "15.22%"
$ws.Range("E26").Value = "4.05%"
$ws.Range("D27").Value = "0.0003395"
$ws.Range("D39").Value = "0.01612"
$ws.Range("D40").Value = "0.04491"
$ws.Range("E40").Value = "2.27%"
$ws.Range("D41").Value = "0.007237"
$ws.Range("E41").Value = "-1.47%"
$ws.Range("D42").Value = "0.009000"
$ws.Range("E43").Value = "1.73%"
$ws.Range("D44").Value = "0.001970"
$ws.Range("E44").Value = "-1.73%"
$ws.Range("D45").Value = "0.009109"
$ws.Range("E45").Value = "-1.33%"
$ws.Range("D46").Value = "0.00006098"
$ws.Range("E46").Value = "2.55%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.06%"
$ws.Range("D48").Value = "2.216"
$ws.Range("E48").Value = "-1.38%"
$ws.Range("D49").Value = "0.002001"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "0.06%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "0.06%"
# Restore the default "Normal" style so the temporary text-format styling
# does not linger on these cells (keeps formatting identical to before).
$priceVolRange.Style = "Normal"
